# Auto-generated script to apply the cryptos.xlsx price/volume update
# described in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force Excel to store the value as literal text (matches the source
    # workbook, where these numbers are stored as inline strings, not
    # numeric cells) and then restore the default cell style so no visible
    # formatting changes are introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "29.956.91"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.893.63"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "0.7773"
$ws.Range("E5").Value = "  +0.13%  "
Set-TextValue $ws.Range("D6") "243.92"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws.Range("D8") "0.3130"
$ws.Range("E8").Value = "  +0.15%  "
Set-TextValue $ws.Range("D9") "25.86"
$ws.Range("E9").Value = "  +1.85%  "
Set-TextValue $ws.Range("D10") "0.07251"
$ws.Range("E10").Value = "  +0.60%  "
Set-TextValue $ws.Range("D11") "0.08709"
$ws.Range("E11").Value = "  +7.80%  "
$ws.Range("D12").Value = "2.135.52"
$ws.Range("E12").Value = "  +12.02%  "
Set-TextValue $ws.Range("D13") "0.7740"
$ws.Range("E13").Value = "  +0.90%  "
Set-TextValue $ws.Range("D14") "5.419"
$ws.Range("E14").Value = "  -1.27%  "
Set-TextValue $ws.Range("D15") "94.56"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").Value = "30.355.60"
$ws.Range("E16").Value = "  +1.72%  "
Set-TextValue $ws.Range("D17") "6.178"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.366.68"
$ws.Range("E18").Value = "  +10.05%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "13.91"
$ws.Range("E19").Value = "  -0.34%  "
Set-TextValue $ws.Range("D20") "246.22"
$ws.Range("E20").Value = "  +1.07%  "
Set-TextValue $ws.Range("D21") "0.000007863"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D22") "1.000"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D23") "8.124"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.02%  "
Set-TextValue $ws.Range("D25") "0.1668"
$ws.Range("E25").Value = "  +6.91%  "
Set-TextValue $ws.Range("D26") "9.486"
$ws.Range("E26").Value = "  +0.97%  "
Set-TextValue $ws.Range("D27") "163.41"
$ws.Range("E27").Value = "  +0.54%  "
Set-TextValue $ws.Range("D28") "18.84"
$ws.Range("E28").Value = "  +0.51%  "
Set-TextValue $ws.Range("D29") "2.054"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  -0.47%  "
Set-TextValue $ws.Range("D32") "4.514"
$ws.Range("E32").Value = "  +0.88%  "
Set-TextValue $ws.Range("D33") "4.139"
$ws.Range("E33").Value = "  +0.97%  "
Set-TextValue $ws.Range("D34") "0.05475"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  -1.31%  "
Set-TextValue $ws.Range("D36") "0.7534"
$ws.Range("E36").Value = "  +0.63%  "
Set-TextValue $ws.Range("D37") "1.009"
$ws.Range("E37").Value = "  +0.90%  "
Set-TextValue $ws.Range("D38") "2.705"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("E40").Value = "  +0.29%  "
Set-TextValue $ws.Range("D41") "0.4512"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "1.112.91"
$ws.Range("E42").Value = "  -2.30%  "
Set-TextValue $ws.Range("D43") "73.61"
$ws.Range("E43").Value = "  +0.01%  "
Set-TextValue $ws.Range("D44") "6.122"
$ws.Range("E44").Value = "  +4.06%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "2.243.02"
$ws.Range("E45").Value = "  +9.60%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D46") "0.8511"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +0.01%  "
Set-TextValue $ws.Range("D48") "103.90"
$ws.Range("E48").Value = "  +0.14%  "
Set-TextValue $ws.Range("D49") "1.877"
$ws.Range("E49").Value = "  -0.73%  "
Set-TextValue $ws.Range("D50") "7.620"
$ws.Range("E50").Value = "  +2.10%  "
Set-TextValue $ws.Range("D51") "9.838"
$ws.Range("E51").Value = "  -0.86%  "
